# Themepark_Followers.xlsx -- "Add files via upload" re-save.
#
# Real content deltas (everything else in the source diff is Excel-generated
# noise: shared-string table reshuffle that leaves every label's text
# unchanged, font-table / cellXf reordering that leaves every cell's actual
# font unchanged, random GUIDs, and window-geometry bookkeeping tied to the
# editing machine):
#   1. B2 (Fun Spot America followers) 66       -> 72
#   2. B60 (canada's wonderland total) 837016   -> 837089
#   3. Trailing empty rows 61:70 removed (used range A1:B70 -> A1:B60)
#   4. Selection/scroll position reset to B2, no frozen topLeftCell
#   5. Column A narrowed, column B given an explicit width
#   6. Workbook theme switched from the newer "Aptos" Office theme back to
#      the classic Calibri-based "Office" theme (colors + major/minor fonts)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) data edits -----------------------------------------------------
$ws.Range("B2").Value = 72
$ws.Range("B60").Value = 837089

# --- 2) drop the trailing blank rows (61-70) ----------------------------
$ws.Range("A61:B70").EntireRow.Delete()

# --- 3) column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.65
$ws.Columns.Item(2).ColumnWidth = 9

# --- 4) selection / scroll position --------------------------------------
$null = $ws.Range("B2").Select()

# --- 5) swap the theme back to the classic Calibri "Office" palette ------
$theme = $wb.Theme
$cs = $theme.ThemeColorScheme
$fs = $theme.ThemeFontScheme

function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# dk1 / lt1 (Colors 1-2) stay as the default black/white system colors.
$cs.Colors(3).RGB  = (RGBVal 0x44 0x54 0x6A)   # dk2
$cs.Colors(4).RGB  = (RGBVal 0xE7 0xE6 0xE6)   # lt2
$cs.Colors(5).RGB  = (RGBVal 0x44 0x72 0xC4)   # accent1
$cs.Colors(6).RGB  = (RGBVal 0xED 0x7D 0x31)   # accent2
$cs.Colors(7).RGB  = (RGBVal 0xA5 0xA5 0xA5)   # accent3
$cs.Colors(8).RGB  = (RGBVal 0xFF 0xC0 0x00)   # accent4
$cs.Colors(9).RGB  = (RGBVal 0x5B 0x9B 0xD5)   # accent5
$cs.Colors(10).RGB = (RGBVal 0x70 0xAD 0x47)   # accent6
$cs.Colors(11).RGB = (RGBVal 0x05 0x63 0xC1)   # hlink
$cs.Colors(12).RGB = (RGBVal 0x95 0x4F 0x72)   # folHlink

$fs.MajorFont(1).Name = "Calibri Light"
$fs.MinorFont(1).Name = "Calibri"

$theme.Name = "Office"
